$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the original content of each data row (2-25) before any writes,
# since several rows are both a source and a destination in the permutation.
$snapshot = @{}
$snapshot[2] = $ws.Range("A2:R2").Value2
$snapshot[3] = $ws.Range("A3:R3").Value2
$snapshot[4] = $ws.Range("A4:R4").Value2
$snapshot[5] = $ws.Range("A5:R5").Value2
$snapshot[6] = $ws.Range("A6:R6").Value2
$snapshot[7] = $ws.Range("A7:R7").Value2
$snapshot[8] = $ws.Range("A8:R8").Value2
$snapshot[9] = $ws.Range("A9:R9").Value2
$snapshot[10] = $ws.Range("A10:R10").Value2
$snapshot[11] = $ws.Range("A11:R11").Value2
$snapshot[12] = $ws.Range("A12:R12").Value2
$snapshot[13] = $ws.Range("A13:R13").Value2
$snapshot[14] = $ws.Range("A14:R14").Value2
$snapshot[15] = $ws.Range("A15:R15").Value2
$snapshot[16] = $ws.Range("A16:R16").Value2
$snapshot[17] = $ws.Range("A17:R17").Value2
$snapshot[18] = $ws.Range("A18:R18").Value2
$snapshot[19] = $ws.Range("A19:R19").Value2
$snapshot[20] = $ws.Range("A20:R20").Value2
$snapshot[21] = $ws.Range("A21:R21").Value2
$snapshot[22] = $ws.Range("A22:R22").Value2
$snapshot[23] = $ws.Range("A23:R23").Value2
$snapshot[24] = $ws.Range("A24:R24").Value2
$snapshot[25] = $ws.Range("A25:R25").Value2

# Apply the permutation: weekly refresh reshuffled the rows into new positions.
$ws.Range("A2:R2").Value2 = $snapshot[15]
$ws.Range("A3:R3").Value2 = $snapshot[22]
$ws.Range("A4:R4").Value2 = $snapshot[18]
$ws.Range("A5:R5").Value2 = $snapshot[19]
$ws.Range("A6:R6").Value2 = $snapshot[20]
$ws.Range("A7:R7").Value2 = $snapshot[17]
$ws.Range("A8:R8").Value2 = $snapshot[24]
$ws.Range("A9:R9").Value2 = $snapshot[25]
$ws.Range("A10:R10").Value2 = $snapshot[2]
$ws.Range("A11:R11").Value2 = $snapshot[9]
$ws.Range("A12:R12").Value2 = $snapshot[16]
$ws.Range("A13:R13").Value2 = $snapshot[21]
$ws.Range("A14:R14").Value2 = $snapshot[11]
$ws.Range("A15:R15").Value2 = $snapshot[7]
$ws.Range("A16:R16").Value2 = $snapshot[13]
$ws.Range("A17:R17").Value2 = $snapshot[14]
$ws.Range("A18:R18").Value2 = $snapshot[10]
$ws.Range("A19:R19").Value2 = $snapshot[23]
$ws.Range("A20:R20").Value2 = $snapshot[5]
$ws.Range("A21:R21").Value2 = $snapshot[8]
$ws.Range("A22:R22").Value2 = $snapshot[12]
$ws.Range("A23:R23").Value2 = $snapshot[3]
$ws.Range("A24:R24").Value2 = $snapshot[6]
$ws.Range("A25:R25").Value2 = $snapshot[4]
